# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# OFF sheet - Home row (row 2) target depth stats updated
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 235
$wsOff.Range("C2").Value = 173
$wsOff.Range("D2").Value = 42
$wsOff.Range("E2").Value = 13
$wsOff.Range("G2").Value = 2

# DEF sheet - Home row (row 2) target depth stats updated
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 196
$wsDef.Range("C2").Value = 143
$wsDef.Range("D2").Value = 36
$wsDef.Range("E2").Value = 19
